$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.314.15"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").Value = "2.265.23"
$ws.Range("E3").Value = "  -0.27%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'317.89"
$ws.Range("E5").Value = "  -1.45%  "
$ws.Range("D6").Value = "'100.08"
$ws.Range("E6").Value = "  -4.61%  "
$ws.Range("E7").Value = "  -2.08%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "'0.540"
$ws.Range("E9").Value = "  -5.38%  "
$ws.Range("D10").Value = "'36.40"
$ws.Range("E10").Value = "  -5.81%  "
$ws.Range("D11").Value = "'0.0828"
$ws.Range("E11").Value = "  -1.79%  "
$ws.Range("D12").Value = "'7.46"
$ws.Range("E12").Value = "  -5.45%  "
$ws.Range("E13").Value = "  -2.36%  "
$ws.Range("D14").Value = "2.610.04"
$ws.Range("E14").Value = "  -0.43%  "
$ws.Range("E15").Value = "  -3.73%  "
$ws.Range("D16").Value = "2.264.14"
$ws.Range("E16").Value = "  -0.43%  "
$ws.Range("E17").Value = "  -3.35%  "
$ws.Range("D18").Value = "44.205.33"
$ws.Range("E18").Value = "  -0.22%  "
$ws.Range("D19").Value = "'13.23"
$ws.Range("E19").Value = "  -4.48%  "
$ws.Range("D20").Value = "0.0₃0983"
$ws.Range("E20").Value = "  -1.86%  "
$ws.Range("D21").Value = "'6.41"
$ws.Range("E21").Value = "  -1.74%  "
$ws.Range("D22").Value = "'65.88"
$ws.Range("E22").Value = "  -0.72%  "
$ws.Range("D23").Value = "'240.28"
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("E24").Value = "  -5.87%  "
$ws.Range("E25").Value = "  -7.52%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.30%  "
$ws.Range("E27").Value = "  +0.25%  "
$ws.Range("D28").Value = "'39.32"
$ws.Range("E28").Value = "  +2.42%  "
$ws.Range("D29").Value = "'2.12"
$ws.Range("E29").Value = "  -4.04%  "
$ws.Range("D30").Value = "'6.11"
$ws.Range("E30").Value = "  -5.79%  "
$ws.Range("E31").Value = "  -2.29%  "
$ws.Range("D32").Value = "'155.23"
$ws.Range("E32").Value = "  -4.48%  "
$ws.Range("D33").Value = "'0.0847"
$ws.Range("E33").Value = "  -3.94%  "
$ws.Range("D34").Value = "'3.50"
$ws.Range("E34").Value = "  +11.90%  "
$ws.Range("D35").Value = "'2.68"
$ws.Range("E35").Value = "  -3.45%  "
$ws.Range("E36").Value = "  -5.38%  "
$ws.Range("D37").Value = "'1.92"
$ws.Range("E37").Value = "  -4.19%  "
$ws.Range("E38").Value = "  -2.18%  "
$ws.Range("D39").Value = "'15.57"
$ws.Range("E39").Value = "  +0.25%  "
$ws.Range("D40").Value = "'3.55"
$ws.Range("E40").Value = "  -9.59%  "
$ws.Range("D41").Value = "'4.01"
$ws.Range("E41").Value = "  -9.27%  "
$ws.Range("E42").Value = "  -5.73%  "
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("D44").Value = "1.737.72"
$ws.Range("E44").Value = "  -2.81%  "
$ws.Range("D45").Value = "'84.27"
$ws.Range("E45").Value = "  -2.13%  "
$ws.Range("E46").Value = "  -5.67%  "
$ws.Range("E47").Value = "  -3.86%  "
$ws.Range("E48").Value = "  -1.27%  "
$ws.Range("D49").Value = "'57.30"
$ws.Range("E49").Value = "  -5.50%  "
$ws.Range("B50").Value = "ordi"
$ws.Range("C50").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D50").Value = "'71.57"
$ws.Range("E50").Value = "  -5.29%  "
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").Value = "'1.63"
$ws.Range("E51").Value = "  -5.26%  "
